# Implement new model to minimize the number of trips to go to the university.
# This works only if you allow more than 1 shift per day.
# Improve input parsing with argparse library.
#
# Re-run the shift-assignment solver output: rename the sheet/title to mark
# this as the "Problem" input, add a new student (Utente Prova 2), and
# rewrite the resulting December 2018 shift assignment grid + shift counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename worksheet tab (workbook.xml <sheet name=.../>)
$ws.Name = "December 2018 Problem"

# Update the title cell to match the new sheet name
$ws.Range("A1").Value = "Automatic assignment for December 2018 Problem"

# --- Week of Dec 3-7 (rows 4-6: 09:30 / 12:30 / 15:30) ---
$ws.Range("C4").Value = "Irene"
$ws.Range("D4").Value = "Utente Prova"
$ws.Range("F4").Value = "Irene"
$ws.Range("C5").Value = "Filippo M"
$ws.Range("E5").Value = "Luigi Berducci"
$ws.Range("C6").Value = "Utente Prova 2"

# --- Week of Dec 10-14 (rows 8-10) ---
$ws.Range("E8").Value = "Utente Prova"
$ws.Range("F8").Value = "Piccola Ketty"
$ws.Range("D9").Value = "Alessandro"
$ws.Range("F9").Value = "Libianchi Gabriele"
$ws.Range("B10").Value = "Luigi Berducci"
$ws.Range("C10").Value = "Irene"
$ws.Range("D10").Value = "Luigi Berducci"

# --- Week of Dec 17-21 (rows 12-14) ---
$ws.Range("B12").Value = "Utente Prova"
$ws.Range("C12").Value = "Andrea Coletta"
$ws.Range("D12").Value = "Alessandro"
$ws.Range("F12").Value = "Agostina"
$ws.Range("B13").Value = "Piccola Ketty"
$ws.Range("C13").Value = "Filippo M"
$ws.Range("D13").Value = "Andrea Coletta"
$ws.Range("F13").Value = "Utente Prova"
$ws.Range("B14").Value = "Utente Prova 2"
$ws.Range("C14").Value = "Utente Prova 2"
$ws.Range("D14").Value = "Libianchi Gabriele"

# --- Week of Dec 24-28 (rows 16-18) ---
$ws.Range("B16").Value = "Agostina"
$ws.Range("C16").Value = "Alessandro"
$ws.Range("D16").Value = "Alessandro"
$ws.Range("E16").Value = "Agostina"
$ws.Range("F16").Value = "Alessandro"
$ws.Range("B17").Value = "Piccola Ketty"
$ws.Range("C17").Value = "Andrea Coletta"
$ws.Range("B18").Value = "Utente Prova 2"
$ws.Range("C18").Value = "Utente Prova 2"
$ws.Range("D18").Value = "Filippo M"
$ws.Range("E18").Value = "Filippo M"

# --- Dec 31 (rows 20-22) ---
$ws.Range("B20").Value = "Agostina"
$ws.Range("B22").Value = "Luigi Berducci"

# --- Shift-count summary table (rows 26-36) ---
$ws.Range("B27").Value = 1
$ws.Range("A28").Value = "Libianchi Gabriele"
$ws.Range("A29").Value = "Irene"
$ws.Range("B29").Value = 5
$ws.Range("A30").Value = "Filippo M"
$ws.Range("A31").Value = "Utente Prova 2"
$ws.Range("B31").Value = 5
$ws.Range("A32").Value = "Utente Prova"
$ws.Range("A33").Value = "Andrea Coletta"
$ws.Range("A34").Value = "Luigi Berducci"
$ws.Range("A35").Value = "Piccola Ketty"

# New student row added to the summary table
$ws.Range("A36").Value = "Alessandro"
$ws.Range("B36").Value = 6
